$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object 'object[,]' 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.049185476250953
$dataBF[0,2] = 1.056631418948794
$dataBF[0,3] = 1.056424711825644
$dataBF[0,4] = 1.06713958029171
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.050042291345222
$dataBF[1,2] = 1.05730876261192
$dataBF[1,3] = 1.057169986226655
$dataBF[1,4] = 1.067924808290429
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.050597475527166
$dataBF[2,2] = 1.05774772046027
$dataBF[2,3] = 1.057653250772435
$dataBF[2,4] = 1.068433933899662
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.050831056983224
$dataBF[3,2] = 1.057932417176784
$dataBF[3,3] = 1.057856657988262
$dataBF[3,4] = 1.068648215034501
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.050870286961656
$dataBF[4,2] = 1.057963437828636
$dataBF[4,3] = 1.057890825137249
$dataBF[4,4] = 1.068684208072889
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.050600595940925
$dataBF[5,2] = 1.05775018776504
$dataBF[5,3] = 1.057655967757966
$dataBF[5,4] = 1.068436796175563
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.049474880967545
$dataBF[6,2] = 1.056860190320937
$dataBF[6,3] = 1.056676368003242
$dataBF[6,4] = 1.067404737180445
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.047497181886414
$dataBF[7,2] = 1.055297123603673
$dataBF[7,3] = 1.05495810116187
$dataBF[7,4] = 1.065594098735263
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.04618282205697
$dataBF[8,2] = 1.054258701051068
$dataBF[8,3] = 1.053818024382454
$dataBF[8,4] = 1.064392491979945
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.045614685333821
$dataBF[9,2] = 1.053809935113369
$dataBF[9,3] = 1.053325671697133
$dataBF[9,4] = 1.063873509350492
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.04540380450102
$dataBF[10,2] = 1.05364337725666
$dataBF[10,3] = 1.053142988479437
$dataBF[10,4] = 1.063680936618225
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.045449032319647
$dataBF[11,2] = 1.053679098397337
$dataBF[11,3] = 1.053182165636036
$dataBF[11,4] = 1.063722235015883
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.045597250779321
$dataBF[12,2] = 1.053796164648624
$dataBF[12,3] = 1.05331056698483
$dataBF[12,4] = 1.063857587125859
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.045688593051441
$dataBF[13,2] = 1.053868310815575
$dataBF[13,3] = 1.053389705654477
$dataBF[13,4] = 1.063941008660825
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.046220548421869
$dataBF[14,2] = 1.054288502826716
$dataBF[14,3] = 1.053850727936394
$dataBF[14,4] = 1.064426963155744
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.046554496186544
$dataBF[15,2] = 1.054552314550763
$dataBF[15,3] = 1.054140266633813
$dataBF[15,4] = 1.064732144698516
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.046749377647782
$dataBF[16,2] = 1.054706276041267
$dataBF[16,3] = 1.054309275736673
$dataBF[16,4] = 1.064910279251151
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.046815843315027
$dataBF[17,2] = 1.054758787218334
$dataBF[17,3] = 1.054366924798533
$dataBF[17,4] = 1.064971040063677
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.04651865687367
$dataBF[18,2] = 1.054524001287733
$dataBF[18,3] = 1.054109188816593
$dataBF[18,4] = 1.064699388422856
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.045553599970232
$dataBF[19,2] = 1.053761687846787
$dataBF[19,3] = 1.053272750492802
$dataBF[19,4] = 1.063817723751093
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.044947701671428
$dataBF[20,2] = 1.053283166160862
$dataBF[20,3] = 1.052747997904439
$dataBF[20,4] = 1.063264547365341
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.045268816653365
$dataBF[21,2] = 1.053536765461753
$dataBF[21,3] = 1.053026069575469
$dataBF[21,4] = 1.063557685880895
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.046534850821979
$dataBF[22,2] = 1.054536794571493
$dataBF[22,3] = 1.054123231155161
$dataBF[22,4] = 1.064714189180716
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.048007748290789
$dataBF[23,2] = 1.055700583522109
$dataBF[23,3] = 1.055401364820919
$dataBF[23,4] = 1.066061234043676
$ws.Range("B2:F25").Value = $dataBF

$dataIM = New-Object 'object[,]' 24,5
$dataIM[0,0] = 1.047076762066251
$dataIM[0,1] = 1.054225262611738
$dataIM[0,2] = 1.059368599651256
$dataIM[0,3] = 1.059162459640861
$dataIM[0,4] = 1.069848244470518
$dataIM[1,0] = 1.047287393483905
$dataIM[1,1] = 1.054731710113886
$dataIM[1,2] = 1.059859967496469
$dataIM[1,3] = 1.059721544458102
$dataIM[1,4] = 1.070449276200888
$dataIM[2,0] = 1.047422918653767
$dataIM[2,1] = 1.05505950156958
$dataIM[2,2] = 1.060177903540013
$dataIM[2,3] = 1.060083662331436
$dataIM[2,4] = 1.070838551887949
$dataIM[3,0] = 1.047479708944336
$dataIM[3,1] = 1.055197324348303
$dataIM[3,2] = 1.060311559632247
$dataIM[3,3] = 1.060235979608428
$dataIM[3,4] = 1.071002289862425
$dataIM[4,0] = 1.047489233443219
$dataIM[4,1] = 1.0552204664889
$dataIM[4,2] = 1.060334000790259
$dataIM[4,3] = 1.060261559168707
$dataIM[4,4] = 1.071029787212795
$dataIM[5,0] = 1.047423678214283
$dataIM[5,1] = 1.055061343090345
$dataIM[5,2] = 1.060179689477588
$dataIM[5,3] = 1.060085697278388
$dataIM[5,4] = 1.070840739425015
$dataIM[6,0] = 1.047148104384165
$dataIM[6,1] = 1.054396400284931
$dataIM[6,2] = 1.059534661322374
$dataIM[6,3] = 1.05935133106428
$dataIM[6,4] = 1.070051288795671
$dataIM[7,0] = 1.046656664256108
$dataIM[7,1] = 1.053225405002875
$dataIM[7,2] = 1.058398010516439
$dataIM[7,3] = 1.058060050600894
$dataIM[7,4] = 1.068663065641188
$dataIM[8,0] = 1.046325154177576
$dataIM[8,1] = 1.052445305568422
$dataIM[8,2] = 1.057640307237373
$dataIM[8,3] = 1.057201146295877
$dataIM[8,4] = 1.067739623045551
$dataIM[9,0] = 1.046180695238632
$dataIM[9,1] = 1.052107664754305
$dataIM[9,2] = 1.057312246121028
$dataIM[9,3] = 1.056829712473978
$dataIM[9,4] = 1.067340266390971
$dataIM[10,0] = 1.046126900373312
$dataIM[10,1] = 1.051982273441457
$dataIM[10,2] = 1.057190395453164
$dataIM[10,3] = 1.056691818667107
$dataIM[10,4] = 1.067192004334722
$dataIM[11,0] = 1.046138445713822
$dataIM[11,1] = 1.05200916921734
$dataIM[11,2] = 1.057216532548081
$dataIM[11,3] = 1.056721394017416
$dataIM[11,4] = 1.067223803547752
$dataIM[12,0] = 1.046176251316014
$dataIM[12,1] = 1.052097299381531
$dataIM[12,2] = 1.057302173779509
$dataIM[12,3] = 1.056818312625721
$dataIM[12,4] = 1.067328009427157
$dataIM[13,0] = 1.046199526546149
$dataIM[13,1] = 1.052151602442
$dataIM[13,2] = 1.057354940984743
$dataIM[13,3] = 1.056878037131384
$dataIM[13,4] = 1.067392224333155
$dataIM[14,0] = 1.046334722273577
$dataIM[14,1] = 1.052467716899966
$dataIM[14,2] = 1.05766208031393
$dataIM[14,3] = 1.057225807324194
$dataIM[14,4] = 1.067766137708895
$dataIM[15,0] = 1.046419283116915
$dataIM[15,1] = 1.052666047578747
$dataIM[15,2] = 1.057854749575137
$dataIM[15,3] = 1.057444083272555
$dataIM[15,4] = 1.068000818745465
$dataIM[16,0] = 1.046468517894188
$dataIM[16,1] = 1.052781744636796
$dataIM[16,2] = 1.057967133037247
$dataIM[16,3] = 1.057571445901033
$dataIM[16,4] = 1.068137752363404
$dataIM[17,0] = 1.046485290708862
$dataIM[17,1] = 1.052821196698332
$dataIM[17,2] = 1.05800545331543
$dataIM[17,3] = 1.057614881015273
$dataIM[17,4] = 1.068184451303439
$dataIM[18,0] = 1.046410219649073
$dataIM[18,1] = 1.052644767102589
$dataIM[18,2] = 1.05783407769298
$dataIM[18,3] = 1.05742065956409
$dataIM[18,4] = 1.067975634699982
$dataIM[19,0] = 1.046165122267122
$dataIM[19,1] = 1.052071346584476
$dataIM[19,2] = 1.057276954398851
$dataIM[19,3] = 1.056789770476813
$dataIM[19,4] = 1.067297321259961
$dataIM[20,0] = 1.046010231124711
$dataIM[20,1] = 1.051710951100655
$dataIM[20,2] = 1.056926703349914
$dataIM[20,3] = 1.056393530040968
$dataIM[20,4] = 1.066871283760793
$dataIM[21,0] = 1.046092416351782
$dataIM[21,1] = 1.051901990143271
$dataIM[21,2] = 1.057112374315467
$dataIM[21,3] = 1.056603543838571
$dataIM[21,4] = 1.067097091558758
$dataIM[22,0] = 1.046414315312469
$dataIM[22,1] = 1.052654382790278
$dataIM[22,2] = 1.057843418418427
$dataIM[22,3] = 1.05743124358787
$dataIM[22,4] = 1.067987014137572
$dataIM[23,0] = 1.046784401128701
$dataIM[23,1] = 1.053528041842384
$dataIM[23,2] = 1.058691856827516
$dataIM[23,3] = 1.058393540500242
$dataIM[23,4] = 1.069021602199975
$ws.Range("I2:M25").Value = $dataIM
